$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("India Road")
$ws.Rows.Item(47).Insert()
$ws.Cells.Item(47, 1).Value = "taxis"
$ws.Cells.Item(47, 2).Formula = '=SUMIFS($H$9:$H$39,$I$9:$I$39,"Taxi")*(''2018 Calcs''!N9/''2018 Calcs''!L9)'
Write-Host "done"
